$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$headers = @(
    "timestamp",
    "submissionid",
    "feedbackType",
    "query",
    "email",
    "phone",
    "firstName",
    "lastName",
    "feedbackText",
    "needsClarification",
    "urgency",
    "impactScope",
    "forwardToDepartment",
    "linkToAdditForm",
    "reminderSent",
    "status"
)
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}
$ws.Range("A1:P1").Font.Bold = $true

# --- Row 2 ---
$ws.Range("A2").Value = 45778.68244238426
$ws.Range("A2").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("B2").Value = "be75975c-2697-11f0-a674-fa163ee583d0"
$ws.Range("E2").Value = "lorism@gmx.net"
$ws.Range("F2").Value = "(078) 715-3999"
$ws.Range("G2").Value = "Julie"
$ws.Range("H2").Value = "Eckmann"
$ws.Range("I2").Value = "Dini mueter"

# --- Row 3 ---
$ws.Range("A3").Value = 45778.73738229166
$ws.Range("A3").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("B3").Value = "cb9b5f8d-26a2-11f0-a674-fa163ee583d0"
$ws.Range("E3").Value = "loris.marino@students.fhnw.ch"
$ws.Range("F3").Value = "(078) 715-3999"
$ws.Range("G3").Value = "Ben"
$ws.Range("H3").Value = "Dover"
$ws.Range("I3").Value = "Test"

# --- Row 4 ---
$ws.Range("A4").Value = 45778.7507708024
$ws.Range("A4").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("B4").Value = "7c24a12b-26a5-11f0-a674-fa163ee583d0"
$ws.Range("E4").Value = "loris.marino@students.fhnw.ch"
$ws.Range("F4").Value = "(078) 715-3999"
$ws.Range("G4").Value = "Loris"
$ws.Range("H4").Value = "Mariño"
$ws.Range("I4").Value = "Das ist die Frage"

Write-Output "edit applied"
